# Generate Report for Handback
# The handback transform failed for the "2451ecbd-690d-472a-b333-b74fe6769850"
# file in both the zh-cn and de-de locales: update the Status cells to reflect
# the failure and record the Error Detail message for each locale sheet.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$zhError = "Handback file name: vso1p42s.vys is different with handoff file name: 2451ecbd-690d-472a-b333-b74fe6769850.29c0a2e6410812595b91451f111f61e5b4079aae.zh-cn."
$deError  = "Handback file name: vso1p42s.vys is different with handoff file name: 2451ecbd-690d-472a-b333-b74fe6769850.29c0a2e6410812595b91451f111f61e5b4079aae.de-de."

# --- Overview sheet: row for 2451ecbd...md now shows the failed status too ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = 40 - 5/6

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 40 - 5/6
